$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.378.18"
$ws.Range("E2").Value = "  -2.93%  "
$ws.Range("D3").Value = "2.286.69"
$ws.Range("E3").Value = "  -2.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "495.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.43"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  -1.88%  "
$ws.Range("D9").Value = "2.284.95"
$ws.Range("E9").Value = "  -3.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0943"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.28%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.323"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("E13").Value = "  -3.89%  "
$ws.Range("D14").Value = "2.692.75"
$ws.Range("E14").Value = "  -2.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.52"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "54.208.44"
$ws.Range("E16").Value = "  -3.14%  "
$ws.Range("E17").Value = "  -3.04%  "
$ws.Range("D18").Value = "2.280.00"
$ws.Range("E18").Value = "  -4.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.95"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.06"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "302.25"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.43"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.39%  "
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("E24").Value = "  -2.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.74"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.51%  "
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").Value = "2.393.97"
$ws.Range("E28").Value = "  -2.67%  "
$ws.Range("E29").Value = "  +2.10%  "
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.84"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.84%  "
$ws.Range("E32").Value = "  -2.98%  "
$ws.Range("D33").Value = "0.0₃0684"
$ws.Range("E33").Value = "  -3.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.88"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.872"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.03%  "
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.61"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.374"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "126.55"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.80"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0888"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("E49").Value = "  -2.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "238.75"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("E51").Value = "  +0.04%  "
